$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "fyOYs108"
$ws.Range("B2").Value = 23110716
$ws.Range("C2").Value = "xvifvfo50"
$ws.Range("D2").Value = "hEM8#9t&"
$ws.Range("F2").Value = "POfOYRdX"
$ws.Range("G2").Value = "zdJO"
